# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
#
# Semantic edits:
#   About!C1                 : last-updated date 45366 (2024-03-15) -> 45379 (2024-03-28)
#   RAF-capacity!B24 (hydrogen combustion turbine) : 0.3 -> 1
#   RAF-capacity!B25 (hydrogen combined cycle)      : 0.3 -> 1
#
# Plus the view state left behind by the author: the RAF-capacity sheet
# becomes the active/selected tab (scrolled down, zoomed to 80%) with B25
# selected, instead of RAF-generation being active with B26 selected.

$wb = $excel.ActiveWorkbook

# --- Update the "About" sheet's last-modified date -------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- Update the hydrogen plant RAF-capacity values --------------------------
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# --- Leave the workbook with RAF-capacity as the active/displayed sheet ----
$wsCapacity.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 80
$wsCapacity.Range("B25").Select()
